$wb = $excel.ActiveWorkbook

$sheetNames = @("展览", "全部类型")

foreach ($name in $sheetNames) {
    $ws = $wb.Worksheets.Item($name)
    $ws.Range("F3").Value = 1794
    $ws.Range("F5").Value = 1111
    $ws.Range("F6").Value = 998
}
